#
# edit.ps1 -- apply the "Agregamos líneas a tus trabajos" edit.
#
# Summary of changes performed:
#   1. Merge "1.-" + "Git revert" runs into a single run.
#   2. Merge the "...borrados no des" + "eados." runs into a single run.
#   3. Remove the mid-paragraph "_GoBack" bookmark and merge the two runs
#      around it into a single run (keeping straight quotes).
#   4. Append four new paragraphs after the "Cuando haces un fork..." text:
#        - an empty paragraph
#        - a paragraph (with a lastRenderedPageBreak marker) of new text
#        - a paragraph of dashed text
#        - a paragraph of dashed text that now carries the "_GoBack" bookmark
#   5. Mint the word/footnotes.xml + word/endnotes.xml parts (with just the
#      standard separator/continuationSeparator boilerplate) the same way
#      Word does internally, by adding a footnote and then deleting it.
#

$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# 1) "1.-" / "Git revert" -> single run "1.-Git revert"
# ---------------------------------------------------------------------
$r1 = $d.Paragraphs(3).Range
$t1 = "1.-Git revert"
$r1.Find.Execute($t1, $false, $false, $false, $false, $false, $true, 1, $false, $t1, 2)

# ---------------------------------------------------------------------
# 2) Merge the split "...borrados no des" / "eados." runs
# ---------------------------------------------------------------------
$r2 = $d.Paragraphs(4).Range
$t2 = $r2.Text.TrimEnd([char]13, [char]7)
$r2.Find.Execute($t2, $false, $false, $false, $false, $false, $true, 1, $false, $t2, 2)

# ---------------------------------------------------------------------
# 3) Remove the "_GoBack" bookmark from the middle of the "fork" sentence
#    and merge the two runs around it into one (preserving straight quotes).
# ---------------------------------------------------------------------
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks("_GoBack").Delete()
}

$p10 = $d.Paragraphs(10).Range
$t3 = 'El "fork" es una de las operativas comunes con el trabajo en Git y GitHub. Básicamente '
$p10.Find.Execute($t3, $false, $false, $false, $false, $false, $true, 1, $false, $t3, 2)

# Typing a replacement through Find re-triggers "smart quotes" autocorrect,
# turning the straight quotes into curly ones -- put the straight quotes
# back with a direct (non-typing) range-text assignment so it doesn't
# re-trigger autocorrect, and so the run stays merged.
$p10b = $d.Paragraphs(10).Range
$quoteFix = $p10b.Duplicate
$openCurly = [char]0x201C
$closeCurly = [char]0x201D
$mangled = $openCurly + "fork" + $closeCurly
if ($quoteFix.Find.Execute($mangled)) {
    $quoteFix.Text = '"fork"'
}

# ---------------------------------------------------------------------
# 4) Append the four new paragraphs after paragraph 12
#    ("Cuando haces un fork de un repositorio...")
# ---------------------------------------------------------------------
$lastExisting = $d.Paragraphs(12).Range
$insertionPoint = $d.Range($lastExisting.End, $lastExisting.End)
$insertionPoint.InsertParagraphAfter()
$insertionPoint.InsertParagraphAfter()
$insertionPoint.InsertParagraphAfter()
$insertionPoint.InsertParagraphAfter()

# --- 4a) a genuinely empty paragraph -----------------------------------
$pBlank = $d.Paragraphs(13).Range
$target13 = $d.Range($pBlank.Start, $pBlank.End)
$xmlBlank = '<?xml version="1.0" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p/></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$target13.InsertXML($xmlBlank)

# --- 4b) new paragraph carrying the page-break render hint -------------
$pBreak = $d.Paragraphs(14).Range
$target14 = $d.Range($pBreak.Start, $pBreak.End)
$text14 = "Líneas nuevas agregadas a tu tarea de git para aumentar un 10 a tu calificación y ser mejor desarrolladora."
$xml14 = '<?xml version="1.0" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r><w:lastRenderedPageBreak/><w:t>' + $text14 + '</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$target14.InsertXML($xml14)

# --- 4c) dashed text line -----------------------------------------------
$d.Paragraphs(15).Range.Text = "------------Líneas mega importantes agregadas--------------------"

# --- 4d) final dashed text line + the relocated "_GoBack" bookmark ------
$pLast = $d.Paragraphs(16).Range
$target16 = $d.Range($pLast.Start, $pLast.End)
$text16 = "------------------Somos cracks-------------"
$xml16 = '<?xml version="1.0" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r><w:t>' + $text16 + '</w:t></w:r><w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$target16.InsertXML($xml16)

# ---------------------------------------------------------------------
# 5) Mint word/footnotes.xml + word/endnotes.xml (separators only) the
#    same way Word does: add a footnote, then remove it again, which
#    leaves the two parts in the package with just their boilerplate.
# ---------------------------------------------------------------------
$fnAnchor = $d.Paragraphs(1).Range
$fnAnchor.Collapse(0)
$fn = $fnAnchor.Footnotes.Add($fnAnchor, "", "x")
$fn.Delete()

Write-Host "Edit complete. Paragraph count:" $d.Paragraphs.Count
